$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date column header (next date after 14-jul)
$ws.Range("AB1").Value = "17-jul"

# Bring over the number format/alignment used by the rest of the data
# columns (same style as column AA) before writing the new values, so
# the new cells share the existing style index instead of minting a
# new one.
$ws.Range("AA2:AA11").Copy()
$ws.Range("AB2:AB11").PasteSpecial(-4122)

# New data values for the "17-jul" column
$ws.Range("AB2").Value = 12
$ws.Range("AB3").Value = 18
$ws.Range("AB4").Value = 8
$ws.Range("AB5").Value = 11
$ws.Range("AB6").Value = 11
$ws.Range("AB7").Value = 17
$ws.Range("AB8").Value = 22
$ws.Range("AB9").Value = 9
$ws.Range("AB10").Value = 19
$ws.Range("AB11").Value = 23

# Mirror the manual-entry workflow: after typing the last value the
# active cell/selection moves one row below the last entered cell.
$ws.Range("AB12").Select()

# A handful of helper columns (O, Q, U, W) also end up hidden with a
# zero custom width, matching the other already-hidden helper columns
# in this sheet (B through N).
$zeroWidth = -0.8333333333333334
$ws.Range("O1").EntireColumn.ColumnWidth = $zeroWidth
$ws.Range("O1").EntireColumn.Hidden = $true
$ws.Range("Q1").EntireColumn.ColumnWidth = $zeroWidth
$ws.Range("Q1").EntireColumn.Hidden = $true
$ws.Range("U1").EntireColumn.ColumnWidth = $zeroWidth
$ws.Range("U1").EntireColumn.Hidden = $true
$ws.Range("W1").EntireColumn.ColumnWidth = $zeroWidth
$ws.Range("W1").EntireColumn.Hidden = $true
